$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source added two new weekly price records for "Cebollín" at
# Femacal de La Calera (date 2023-04-25 / serial 45041), one "Primera"
# and one "Segunda" grade row, inserted right before the existing
# row for serial 44705 (2022-05-24). Inserting at row 659 pushes all
# the old rows 659:694 down to 661:696 and grows the sheet dimension
# from A1:R694 to A1:R696, matching the target diff exactly.

$ws.Rows.Item(659).Insert()
$ws.Rows.Item(659).Insert()

# --- New row 659: Cebollín, Primera ---
$ws.Range("A659").Value = 3
$ws.Range("B659").Value = "Femacal de La Calera"
$ws.Range("C659").Value = "Coquimbo"
$ws.Range("D659").Value = 45041
$ws.Range("E659").Value = 5
$ws.Range("F659").Value = 100112037
$ws.Range("G659").Value = "Cebollín"
$ws.Range("H659").Value = "Sin especificar"
$ws.Range("I659").Value = "Primera"
$ws.Range("J659").Value = 245
$ws.Range("K659").Value = 3500
$ws.Range("L659").Value = 4000
$ws.Range("M659").Value = 3745
$ws.Range("N659").Value = "$/paquete 36 unidades"
$ws.Range("O659").Value = "Provincia de Quillota"
$ws.Range("P659").Value = 104
$ws.Range("Q659").Value = 36
$ws.Range("R659").Value = "Hortaliza"

# --- New row 660: Cebollín, Segunda ---
$ws.Range("A660").Value = 3
$ws.Range("B660").Value = "Femacal de La Calera"
$ws.Range("C660").Value = "Coquimbo"
$ws.Range("D660").Value = 45041
$ws.Range("E660").Value = 5
$ws.Range("F660").Value = 100112037
$ws.Range("G660").Value = "Cebollín"
$ws.Range("H660").Value = "Sin especificar"
$ws.Range("I660").Value = "Segunda"
$ws.Range("J660").Value = 120
$ws.Range("K660").Value = 3000
$ws.Range("L660").Value = 3000
$ws.Range("M660").Value = 3000
$ws.Range("N660").Value = "$/paquete 36 unidades"
$ws.Range("O660").Value = "Provincia de Quillota"
$ws.Range("P660").Value = 83
$ws.Range("Q660").Value = 36
$ws.Range("R660").Value = "Hortaliza"
